$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = 'questions = [
    {
        "title": "You are the wealth advisor preparing the portfolio construction process of an endowment that invests in ETFs. You came across different regional ETFs that track the performance of the MSCI Emerging Market Index. The chart below shows the main characteristics of the ETFs\u2019 fact sheets.The trustees told you that they prefer to invest in the ETF with the lowest tracking error. Which one should you choose for the trust?",
        "ques_type": 2,
        "options": [
            "ETF A",
            "ETF B",
            "ETF C",
            "ETF D"
        ],
        "score": "ETF B"
    },
    {
        "title": "Andrew bought some stocks last year at an initial cost of $30,000. The current value of those holdings is $20,000. As his wealth advisor, you advise Andrew that he can sell those securities and use the realized loss to offset other gains and decrease his tax bill. Andrew wants to know how much this strategy will reduce his tax payment, assuming a 15% tax rate. What should you tell him?",
        "ques_type": 2,
        "options": [
            "$8,500",
            "$4,500",
            "$3,000",
            "$1,500"
        ],
        "score": "$1,500"
    },
    {
        "title": "You manage the portfolio of a client who just received his annual bonus of $10,000. The client wants to add real estate exposure to their portfolio. They want to invest in a security that has low unsystematic risk and can be sold when the client needs money.What type of assets should you recommend to the client to buy?",
        "ques_type": 2,
        "options": [
            "Publicly listed stock of a company that manages several real estate properties.",
            "Income-producing real estate (i.e., an apartment that he rents to others).",
            "Bond of a real estate company.",
            "ETF that exclusively holds real estate investments trusts (REITs)."
        ],
        "score": "ETF that exclusively holds real estate investments trusts (REITs)."
    },
    {
        "title": "Your client holds a $1m position in Amazon (AMZN) and has $500,000 invested in Treasury bonds. They want to hold AMZN for the long run but expect that the value of AMZN might decrease over the next six months, while the bonds\u2019 value is expected to increase. The client asks you about the best strategy to follow during the next six months to decrease their AMZN exposure to $500,000 and increase their bond exposure to $1m while minimizing transaction costs. What strategy should you recommend?",
        "ques_type": 2,
        "options": [
            "Enter a swap agreement to pay the return on $500,000 worth of AMZN shares and receive the return on a Treasury bond index for six months.",
            "Buy a call option expiring in six months on $500,000 of Treasury bonds and sell a put option on $500,000 of AMZN.",
            "Sell $500,000 of AMZN shares and buy $500,000 of Treasury bonds for six months.",
            "Enter a forward contract to sell $500,000 of AMZN and buy $500,000 of Treasury bonds."
        ],
        "score": "Enter a swap agreement to pay the return on $500,000 worth of AMZN shares and receive the return on a Treasury bond index for six months."
    }
]'
$ws.Rows.Item(1).AutoFit()
